$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2098.0598
$ws.Range("I15").Value = 2098.0598
$ws.Range("K15").Value = 6294.1794
$ws.Range("M15").Value = -6125.1794
$ws.Range("H87").Value = 121107
$ws.Range("J87").Value = 121107
$ws.Range("L87").Value = 121107
$ws.Range("N87").Value = -123603
$ws.Range("H90").Value = 121107
$ws.Range("J90").Value = 121107
$ws.Range("L90").Value = 363321
$ws.Range("N90").Value = -375801
$ws.Range("H98").Value = 43479116
$ws.Range("I98").Value = 47619884
$ws.Range("K98").Value = 47619884
$ws.Range("M98").Value = -47618386
$ws.Range("H122").Value = 43479116
$ws.Range("I122").Value = 47619884
$ws.Range("K122").Value = 142859652
$ws.Range("M122").Value = -142857202

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1622.5186
$ws.Range("I2").Value = 1660.75
$ws.Range("J2").Value = 1513.2858
$ws.Range("K2").Value = 1660.75
$ws.Range("L2").Value = 1513.2858
$ws.Range("M2").Value = -1547.75
$ws.Range("N2").Value = -1739.2858
$ws.Range("H4").Value = 1467
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 1886.41
$ws.Range("I32").Value = 1886.41
$ws.Range("K32").Value = 1886.41
$ws.Range("M32").Value = -1599.41
$ws.Range("H61").Value = 46881050
$ws.Range("I61").Value = 41670628
$ws.Range("K61").Value = 41670628
$ws.Range("M61").Value = -41670416
$ws.Range("H109").Value = 73492
$ws.Range("J109").Value = 73492
$ws.Range("L109").Value = 73492
$ws.Range("N109").Value = -76266
$ws.Range("H112").Value = 21748.75
$ws.Range("J112").Value = 21748.75
$ws.Range("L112").Value = 21748.75
$ws.Range("N112").Value = -24702.75
$ws.Range("H116").Value = 1622.5186
$ws.Range("I116").Value = 1660.75
$ws.Range("J116").Value = 1513.2858
$ws.Range("K116").Value = 1660.75
$ws.Range("L116").Value = 1513.2858
$ws.Range("M116").Value = 633.25
$ws.Range("N116").Value = -6101.2858
$ws.Range("H122").Value = 2506.8086
$ws.Range("I122").Value = 1984.2941
$ws.Range("K122").Value = 5952.8823
$ws.Range("M122").Value = -3502.8823
$ws.Range("H128").Value = 119499.5
$ws.Range("J128").Value = 119499.5
$ws.Range("L128").Value = 119499.5
$ws.Range("N128").Value = -129459.5
$ws.Range("H131").Value = 82248
$ws.Range("J131").Value = 82248
$ws.Range("L131").Value = 82248
$ws.Range("N131").Value = -92328
$ws.Range("H132").Value = 8776809
$ws.Range("I132").Value = 12348266
$ws.Range("K132").Value = 37044798
$ws.Range("M132").Value = -37042268
$ws.Range("H136").Value = 46881050
$ws.Range("I136").Value = 41670628
$ws.Range("K136").Value = 125011884
$ws.Range("M136").Value = -125009334
$ws.Range("H141").Value = 119500
$ws.Range("J141").Value = 119500
$ws.Range("L141").Value = 119500
$ws.Range("N141").Value = -129860

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1622.5186
$ws.Range("I3").Value = 1660.75
$ws.Range("J3").Value = 1513.2858
$ws.Range("K3").Value = 1660.75
$ws.Range("L3").Value = 1513.2858
$ws.Range("M3").Value = -1546.75
$ws.Range("N3").Value = -1741.2858
$ws.Range("H94").Value = 1138.4286
$ws.Range("I94").Value = 538.1786
$ws.Range("K94").Value = 538.1786
$ws.Range("M94").Value = -87.17859999999996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5523.4
$ws.Range("I86").Value = 5400.2856
$ws.Range("J86").Value = 5631.125
$ws.Range("K86").Value = 5400.2856
$ws.Range("L86").Value = 5631.125
$ws.Range("M86").Value = -4277.2856
$ws.Range("N86").Value = -7877.125
$ws.Range("H89").Value = 5523.4
$ws.Range("I89").Value = 5400.2856
$ws.Range("J89").Value = 5631.125
$ws.Range("K89").Value = 27001.428
$ws.Range("L89").Value = 28155.625
$ws.Range("M89").Value = -21385.428
$ws.Range("N89").Value = -39387.625
$ws.Range("H94").Value = 3387.4583
$ws.Range("I94").Value = 2961.1428
$ws.Range("K94").Value = 2961.1428
$ws.Range("M94").Value = -2510.1428
$ws.Range("H134").Value = 2495.5134
$ws.Range("I134").Value = 1639.826
$ws.Range("J134").Value = 3901.2856
$ws.Range("K134").Value = 4919.478
$ws.Range("L134").Value = 11703.8568
$ws.Range("M134").Value = -2384.478
$ws.Range("N134").Value = -16773.8568
$ws.Range("H141").Value = 173259.38
$ws.Range("J141").Value = 199472.2
$ws.Range("L141").Value = 199472.2
$ws.Range("N141").Value = -209832.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8535345
$ws.Range("I4").Value = 6958730
$ws.Range("J4").Value = 11951344
$ws.Range("K4").Value = 20876190
$ws.Range("L4").Value = 35854032
$ws.Range("M4").Value = -20876078
$ws.Range("N4").Value = -35854256
$ws.Range("H5").Value = 2062.5293
$ws.Range("I5").Value = 1761.75
$ws.Range("J5").Value = 2155.077
$ws.Range("K5").Value = 5285.25
$ws.Range("L5").Value = 6465.231000000001
$ws.Range("M5").Value = -5173.25
$ws.Range("N5").Value = -6689.231000000001
$ws.Range("H56").Value = 6382.857
$ws.Range("I56").Value = 6382.857
$ws.Range("K56").Value = 6382.857
$ws.Range("M56").Value = -5852.857
$ws.Range("H134").Value = 11363.036
$ws.Range("I134").Value = 10674.714
$ws.Range("J134").Value = 13428
$ws.Range("K134").Value = 32024.142
$ws.Range("L134").Value = 40284
$ws.Range("M134").Value = -26954.142
$ws.Range("N134").Value = -50424
$ws.Range("H135").Value = 2062.5293
$ws.Range("I135").Value = 1761.75
$ws.Range("J135").Value = 2155.077
$ws.Range("K135").Value = 15855.75
$ws.Range("L135").Value = 19395.693
$ws.Range("M135").Value = -13320.75
$ws.Range("N135").Value = -24465.693
$ws.Range("H140").Value = 99181.45
$ws.Range("I140").Value = 99181.45
$ws.Range("K140").Value = 297544.35
$ws.Range("M140").Value = -292364.35
$ws.Range("H141").Value = 439282.84
$ws.Range("I141").Value = 757495
$ws.Range("K141").Value = 2272485
$ws.Range("M141").Value = -2267305

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 215625
$ws.Range("I19").Value = 286666.66
$ws.Range("J19").Value = 2500
$ws.Range("K19").Value = 286666.66
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = -286378.66
$ws.Range("N19").Value = -3076
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H113").Value = 4682.615
$ws.Range("I113").Value = 4455.25
$ws.Range("K113").Value = 4455.25
$ws.Range("M113").Value = -2285.25
$ws.Range("H132").Value = 83337790
$ws.Range("I132").Value = 83337790
$ws.Range("K132").Value = 250013370
$ws.Range("M132").Value = -250010840

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 53998
$ws.Range("J43").Value = 54994
$ws.Range("L43").Value = 54994
$ws.Range("N43").Value = -55380
$ws.Range("H82").Value = 1742.875
$ws.Range("I82").Value = 754
$ws.Range("J82").Value = 3014.2856
$ws.Range("K82").Value = 754
$ws.Range("L82").Value = 3014.2856
$ws.Range("M82").Value = -393
$ws.Range("N82").Value = -3736.2856
$ws.Range("H85").Value = 1742.875
$ws.Range("I85").Value = 754
$ws.Range("J85").Value = 3014.2856
$ws.Range("K85").Value = 754
$ws.Range("L85").Value = 3014.2856
$ws.Range("M85").Value = 494
$ws.Range("N85").Value = -5510.2856
$ws.Range("H98").Value = 30177.5
$ws.Range("J98").Value = 30177.5
$ws.Range("L98").Value = 30177.5
$ws.Range("N98").Value = -36167.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 633738.4399999999
$ws.Range("I132").Value = 8930.9
$ws.Range("K132").Value = 26792.7
$ws.Range("M132").Value = -24262.7
